# Update automate_finance.qmd output SOR workbook to use refreshed CVD files.
# This mirrors the refreshed "ytd" (column E) figures and the refreshed
# monthly split values pulled in by the new CVD source data.

$wb = $excel.ActiveWorkbook

# --- Sheet: Baja California Mexico -----------------------------------
$ws = $wb.Worksheets.Item("Baja California Mexico")

# Row 4 (Manufacturing Voluntary Turnover / AOP) - refreshed monthly split
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = $null
$ws.Range("J4").Value = 0.6667
$ws.Range("K4").Value = $null
$ws.Range("L4").Value = $null

# Row 7-9 ytd (column E) refreshed
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776

# Row 9 (Manufacturing Voluntary Turnover / Commit-Forecast) monthly split zeroed out
$ws.Range("I9").Value = $null
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = $null
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 0
$ws.Range("U9").Value = 0
$ws.Range("V9").Value = 0
$ws.Range("W9").Value = 0

# --- Sheet: Bristol Connecticut ---------------------------------------
$ws = $wb.Worksheets.Item("Bristol Connecticut")
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776

# --- Sheet: Charlotte  North Carolina ---------------------------------
$ws = $wb.Worksheets.Item("Charlotte  North Carolina")
$ws.Range("L4").Value = $null

# --- Sheet: Cleveland Ohio ---------------------------------------------
$ws = $wb.Worksheets.Item("Cleveland Ohio")
$ws.Range("E4").Value = 0.0776
$ws.Range("E5").Value = 0.0776
$ws.Range("E6").Value = 0.0776
$ws.Range("L6").Value = $null

# --- Sheet: Marengo Illinois --------------------------------------------
$ws = $wb.Worksheets.Item("Marengo Illinois")
$ws.Range("L4").Value = $null
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("L9").Value = $null

# --- Sheet: Fremont California ------------------------------------------
$ws = $wb.Worksheets.Item("Fremont California")
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776
